# Ajout de cas d'utilisations
#
# Target change (paragraph "Nom : Afficher le stock des ingrédients (PACKAGE : )"):
#   - " (PACKAGE\xa0:"  ->  " (PACKAGE\xa0: "   (trailing space added after the colon)
#   - "  )" (two spaces + close paren) is split into "Gestion" + a relocated
#     "_GoBack" bookmark + " )" (space + close paren)
#   - The document's old "_GoBack" bookmark (previously sitting right after
#     "Ergonomie : ") is removed, since Word keeps only one "_GoBack" bookmark,
#     tracking the most recent edit location.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the stale "_GoBack" bookmark near "Ergonomie :".
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the "(PACKAGE" run inside the "Nom :" paragraph.
# ------------------------------------------------------------------
$nbsp = [char]0xA0
$full = $d.Content.Text
$idx = $full.IndexOf("(PACKAGE")

$runStart = $idx - 1                                   # leading space before "("
$oldRunAText = " (PACKAGE" + $nbsp + ":"
$runAEnd = $runStart + $oldRunAText.Length              # position right after the colon

# Sanity-check we found the expected text before touching anything.
$check = $d.Range($runStart, $runAEnd)
if ($check.Text -ne $oldRunAText) {
    throw "Unexpected document content near '(PACKAGE' run: '$($check.Text)'"
}

# ------------------------------------------------------------------
# 3. Insert a trailing space after the colon: " (PACKAGE: " -> " (PACKAGE : ".
#    A temporary bookmark pins the run's left edge so the edit doesn't get
#    coalesced into the previous run (both share the same run formatting).
# ------------------------------------------------------------------
$boundary = $d.Range($runStart, $runStart)
$d.Bookmarks.Add("ZZTMP_BOUNDARY", $boundary)

$insertPoint = $d.Range($runAEnd, $runAEnd)
$insertPoint.InsertAfter(" ")

$d.Bookmarks("ZZTMP_BOUNDARY").Delete()

# ------------------------------------------------------------------
# 4. The following run used to hold "  )" (two spaces, close paren) right
#    after the just-edited run (now one char further because of the
#    inserted space). Re-seat "_GoBack" in between its two halves — this
#    both splits the run in two (so later edits stay local) and restores
#    the bookmark at its new location.
# ------------------------------------------------------------------
$oldRunBStart = $runAEnd + 1
$splitPoint = $d.Range($oldRunBStart + 1, $oldRunBStart + 1)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# ------------------------------------------------------------------
# 5. Turn the first half (a lone space) into "Gestion".
# ------------------------------------------------------------------
$runB = $d.Range($oldRunBStart, $oldRunBStart + 1)
if ($runB.Text -ne " ") {
    throw "Unexpected document content where 'Gestion' should go: '$($runB.Text)'"
}
$runB.Text = "Gestion"
